# Append: 2025-10-15 18:25 JST
# Adds 3 new rows to the "ランサーズ" (Lancers) listing sheet:
#   - one new row inserted in the middle (after old row 7 / before old row 8)
#   - two new rows appended near the end (before the old last row)
# and refreshes every row's "取得日時" (fetched-at) timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-10-15 18:25:50"

# --- 1. Insert a new data row after the current row 7 (becomes row 8), ---
#        pushing the old rows 8-12 down to rows 9-13.
$ws.Rows(8).Insert()

# --- 2. Insert two more rows just before the (now shifted) old last row ---
#        (old row 12 "Access..." is currently at row 13); this pushes it
#        down to row 15, opening up rows 13 and 14 for new data.
$ws.Rows(13).Insert()
$ws.Rows(13).Insert()

# --- 3. Fill in the brand-new row 8 (コスパスポーツジムの高速自動予約botの開発) ---
$ws.Range("A8").Value = $newTimestamp
$ws.Range("B8").Value = "コスパスポーツジムの高速自動予約botの開発"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5413835"
$ws.Range("G8").Value = 173
$ws.Range("H8").Value = "★bot ◆開発"

# --- 4. Fill in the brand-new row 13 (【急募】16タイプ診断コンテンツのLP制作) ---
$ws.Range("A13").Value = $newTimestamp
$ws.Range("B13").Value = "【急募】16タイプ診断コンテンツのLP制作"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5408735"
$ws.Range("G13").Value = 25

# --- 5. Fill in the brand-new row 14 (【急募】Teams Roomsの設定設置と保守サポート依頼) ---
$ws.Range("A14").Value = $newTimestamp
$ws.Range("B14").Value = "【急募】Teams Roomsの設定設置と保守サポート依頼"
$ws.Range("C14").Value = "システム開発"
$ws.Range("D14").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E14").Value = "期限情報なし"
$ws.Range("F14").Value = "https://www.lancers.jp/work/detail/5408814"
$ws.Range("G14").Value = 18

# --- 6. Refresh the "取得日時" timestamp on every pre-existing data row ---
#        (rows 2-7 were never touched by the row-inserts above; rows 9-12
#        and 15 are the shifted former rows 8-12).
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp
$ws.Range("A9").Value = $newTimestamp
$ws.Range("A10").Value = $newTimestamp
$ws.Range("A11").Value = $newTimestamp
$ws.Range("A12").Value = $newTimestamp
$ws.Range("A15").Value = $newTimestamp

# --- 7. Rebuild the hyperlinks collection from scratch, in row order. ---
#        Row-insert can leave the worksheet's hyperlink relationships out of
#        sync with the shifted cells, so clear everything and re-add each
#        link fresh against its final, correct cell.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5413825")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5413230")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5413280")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5413402")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5413210")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5413215")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5413835")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5413508")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5413293")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5016989")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5413043")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5408735")
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5408814")
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5413333")
